$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the navigation XPath strings for the motorcycle / camper links so
#    they target the "nav_motorcycle" / "nav_camper" element ids instead of
#    the old "motorcycle" / "camper" ids (shared strings used by I1 / J1).
# ---------------------------------------------------------------------------
$ws.Range("I1").Value = '//*[@class="main-navigation"]//*[@id="nav_motorcycle"]'
$ws.Range("J1").Value = '//*[@class="main-navigation"]//*[@id="nav_camper"]'

# Widen columns I and J to fit the new (longer) header text.
$ws.Columns("I").ColumnWidth = 48.98
$ws.Columns("J").ColumnWidth = 41.83

# ---------------------------------------------------------------------------
# 2. Move the current selection from G9 to J7.
# ---------------------------------------------------------------------------
$ws.Range("J7").Select()

# ---------------------------------------------------------------------------
# 3. Reposition the smoke-test screenshot picture: it shifts left slightly
#    and moves further down the sheet, while keeping the same physical size.
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 0
$shp.Top = 161.4
$shp.Width = 1498.9603937007873
$shp.Height = 726.03968503937
